$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 5.28
$ws.Range("C5").Value = 3.63
$ws.Range("D5").Value = 0.73
$ws.Range("E5").Value = 11.18
$ws.Range("F5").Value = 9.26
$ws.Range("G5").Value = 4.16
$ws.Range("H5").Value = 22.26
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 2.77
$ws.Range("K5").Value = 4.09
$ws.Range("L5").Value = 4.58
$ws.Range("M5").Value = 4.67
$ws.Range("N5").Value = 1.34
$ws.Range("O5").Value = 4.14
$ws.Range("P5").Value = 5.86
$ws.Range("Q5").Value = 3.66
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.38
$ws.Range("T5").Value = 55.68
$ws.Range("U5").Value = 11.86
$ws.Range("V5").Value = 3.82
$ws.Range("W5").Value = 7.82
$ws.Range("X5").Value = 4.22
$ws.Range("Y5").Value = 0.44
$ws.Range("Z5").Value = 10.02
$ws.Range("AA5").Value = 3.37
$ws.Range("AB5").Value = 3.11
$ws.Range("AC5").Value = 3.63
$ws.Range("AD5").Value = 4.74
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 20.43
$ws.Range("AG5").Value = 2.06
$ws.Range("AH5").Value = 4.77

# Remove the last data row (row 6) entirely - "1000개" dataset trim
$ws.Rows.Item(6).Delete()
